$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 876.13336
$ws.Range("I58").Value = 253
$ws.Range("J58").Value = 1102.7273
$ws.Range("K58").Value = 759
$ws.Range("L58").Value = 3308.1819
$ws.Range("M58").Value = -609
$ws.Range("N58").Value = -3608.1819

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3379.6316
$ws.Range("I70").Value = 994.2222
$ws.Range("J70").Value = 5526.5
$ws.Range("K70").Value = 2982.6666
$ws.Range("L70").Value = 16579.5
$ws.Range("M70").Value = -2712.6666
$ws.Range("N70").Value = -17119.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3379.6316
$ws.Range("I73").Value = 994.2222
$ws.Range("J73").Value = 5526.5
$ws.Range("K73").Value = 2982.6666
$ws.Range("L73").Value = 16579.5
$ws.Range("M73").Value = -2046.6666
$ws.Range("N73").Value = -18451.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1139.7561
$ws.Range("J127").Value = 1510.8462
$ws.Range("L127").Value = 4532.5386
$ws.Range("N127").Value = -14452.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1015.63464
$ws.Range("J129").Value = 1082.6383
$ws.Range("L129").Value = 3247.9149
$ws.Range("N129").Value = -13247.9149

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4406.9497
$ws.Range("I131").Value = 274.125
$ws.Range("J131").Value = 4770.275
$ws.Range("K131").Value = 822.375
$ws.Range("L131").Value = 14310.825
$ws.Range("M131").Value = 4217.625
$ws.Range("N131").Value = -24390.825

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 26612.96
$ws.Range("I132").Value = 31879.477
$ws.Range("J132").Value = 2035.8889
$ws.Range("K132").Value = 95638.431
$ws.Range("L132").Value = 6107.6667
$ws.Range("M132").Value = -93108.431
$ws.Range("N132").Value = -11167.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1239.8667
$ws.Range("I137").Value = 1094.5435
$ws.Range("J137").Value = 1717.3572
$ws.Range("K137").Value = 3283.6305
$ws.Range("L137").Value = 5152.071599999999
$ws.Range("M137").Value = -733.6305000000002
$ws.Range("N137").Value = -10252.0716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 108790
$ws.Range("J140").Value = 108790
$ws.Range("L140").Value = 108790
$ws.Range("N140").Value = -119150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6279.1724
$ws.Range("I32").Value = 5552
$ws.Range("J32").Value = 13460
$ws.Range("K32").Value = 5552
$ws.Range("L32").Value = 13460
$ws.Range("M32").Value = -5265
$ws.Range("N32").Value = -14034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 868.7
$ws.Range("I97").Value = 742.5925999999999
$ws.Range("J97").Value = 2003.6666
$ws.Range("K97").Value = 742.5925999999999
$ws.Range("L97").Value = 2003.6666
$ws.Range("M97").Value = -246.5925999999999
$ws.Range("N97").Value = -2995.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2389.7585
$ws.Range("I86").Value = 2633.111
$ws.Range("J86").Value = 1991.5454
$ws.Range("K86").Value = 2633.111
$ws.Range("L86").Value = 1991.5454
$ws.Range("M86").Value = -1510.111
$ws.Range("N86").Value = -4237.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2389.7585
$ws.Range("I89").Value = 2633.111
$ws.Range("J89").Value = 1991.5454
$ws.Range("K89").Value = 13165.555
$ws.Range("L89").Value = 9957.726999999999
$ws.Range("M89").Value = -7549.555
$ws.Range("N89").Value = -21189.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 26740.342
$ws.Range("I94").Value = 34796.656
$ws.Range("J94").Value = 781.1111
$ws.Range("K94").Value = 34796.656
$ws.Range("L94").Value = 781.1111
$ws.Range("M94").Value = -34345.656
$ws.Range("N94").Value = -1683.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1029.6666
$ws.Range("I12").Value = 1029.6666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1029.6666
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -859.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2083.5334
$ws.Range("I31").Value = 1056.766
$ws.Range("J31").Value = 3807.0356
$ws.Range("K31").Value = 1056.766
$ws.Range("L31").Value = 3807.0356
$ws.Range("M31").Value = -761.7660000000001
$ws.Range("N31").Value = -4397.0356

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2083.5334
$ws.Range("I34").Value = 1056.766
$ws.Range("J34").Value = 3807.0356
$ws.Range("K34").Value = 1056.766
$ws.Range("L34").Value = 3807.0356
$ws.Range("M34").Value = -854.7660000000001
$ws.Range("N34").Value = -4211.0356

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1115.9103
$ws.Range("I58").Value = 894.38983
$ws.Range("J58").Value = 1803.7894
$ws.Range("K58").Value = 894.38983
$ws.Range("L58").Value = 1803.7894
$ws.Range("M58").Value = -691.38983
$ws.Range("N58").Value = -2209.7894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1760.0435
$ws.Range("I105").Value = 2251.4285
$ws.Range("J105").Value = 995.6667
$ws.Range("K105").Value = 2251.4285
$ws.Range("L105").Value = 995.6667
$ws.Range("M105").Value = -504.4285
$ws.Range("N105").Value = -4489.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1335.3513
$ws.Range("I134").Value = 1249.7627
$ws.Range("J134").Value = 1672
$ws.Range("K134").Value = 3749.2881
$ws.Range("L134").Value = 5016
$ws.Range("M134").Value = -1214.2881
$ws.Range("N134").Value = -10086

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1115.9103
$ws.Range("I136").Value = 894.38983
$ws.Range("J136").Value = 1803.7894
$ws.Range("K136").Value = 2683.16949
$ws.Range("L136").Value = 5411.3682
$ws.Range("M136").Value = -133.1694899999998
$ws.Range("N136").Value = -10511.3682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4775.5454
$ws.Range("I62").Value = 3512
$ws.Range("J62").Value = 5249.375
$ws.Range("K62").Value = 10536
$ws.Range("L62").Value = 15748.125
$ws.Range("M62").Value = -9850
$ws.Range("N62").Value = -17120.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 4775.5454
$ws.Range("I65").Value = 3512
$ws.Range("J65").Value = 5249.375
$ws.Range("K65").Value = 31608
$ws.Range("L65").Value = 47244.375
$ws.Range("M65").Value = -28176
$ws.Range("N65").Value = -54108.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 901494.25
$ws.Range("I107").Value = 683.6667
$ws.Range("J107").Value = 2252710.2
$ws.Range("K107").Value = 2051.0001
$ws.Range("L107").Value = 6758130.600000001
$ws.Range("M107").Value = -131.0001000000002
$ws.Range("N107").Value = -6761970.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1813.174
$ws.Range("I122").Value = 484.92307
$ws.Range("J122").Value = 3539.9
$ws.Range("K122").Value = 4364.30763
$ws.Range("L122").Value = 31859.1
$ws.Range("M122").Value = -1914.30763
$ws.Range("N122").Value = -36759.10000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2349.3333
$ws.Range("I131").Value = 576.5833
$ws.Range("J131").Value = 2593.8506
$ws.Range("K131").Value = 1729.7499
$ws.Range("L131").Value = 7781.551800000001
$ws.Range("M131").Value = 3310.2501
$ws.Range("N131").Value = -17861.5518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1465247.9
$ws.Range("I14").Value = 3346325
$ws.Range("J14").Value = 54440
$ws.Range("K14").Value = 3346325
$ws.Range("L14").Value = 54440
$ws.Range("M14").Value = -3346157
$ws.Range("N14").Value = -54776

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4148.3335
$ws.Range("I62").Value = 4300
$ws.Range("J62").Value = 3996.6667
$ws.Range("K62").Value = 4300
$ws.Range("L62").Value = 3996.6667
$ws.Range("M62").Value = -3676
$ws.Range("N62").Value = -5244.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4148.3335
$ws.Range("I65").Value = 4300
$ws.Range("J65").Value = 3996.6667
$ws.Range("K65").Value = 21500
$ws.Range("L65").Value = 19983.3335
$ws.Range("M65").Value = -18380
$ws.Range("N65").Value = -26223.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 9135.5
$ws.Range("J68").Value = 9135.5
$ws.Range("L68").Value = 9135.5
$ws.Range("N68").Value = -10757.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 9135.5
$ws.Range("J71").Value = 9135.5
$ws.Range("L71").Value = 27406.5
$ws.Range("N71").Value = -35518.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1481.1818
$ws.Range("I81").Value = 1200
$ws.Range("J81").Value = 1543.6666
$ws.Range("K81").Value = 2400
$ws.Range("L81").Value = 3087.3332
$ws.Range("M81").Value = -1339
$ws.Range("N81").Value = -5209.3332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1481.1818
$ws.Range("I84").Value = 1200
$ws.Range("J84").Value = 1543.6666
$ws.Range("K84").Value = 12000
$ws.Range("L84").Value = 15436.666
$ws.Range("M84").Value = -6696
$ws.Range("N84").Value = -26044.666
